$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Logitech G203 Mouse Gaming
$ws.Range("A3").Value = 'Logitech G203 Mouse Gaming'
$desc1 = @'
New Version of Logitech G103 Macro Gaming Mouse
Garansi Perangkat Keras Terbatas 2 Tahun
Tolong simpan dan sertakan kemasan dan nota. Untuk komplain seperti salah barang dan cacat fisik, tolong menyertakan video unboxing.
SIAP UNTUK BERMAIN
Optimalkan waktu bermainmu dengan G203 gaming mouse yang dilengkapi dengan teknologi LIGHTSYNC, sensor kelas gaming, dan desain klasik 6 tombol. Ceriakan game-mu … dan mejamu.
LIGHTSYNC RGB COLOR WAVE
Pilihlah dari pengaturan berdasarkan game dan media, serta animasi yang ceria, atau programlah pengaturanmu sendiri dari sekitar 16,8 juta warna.
KUSTOMISASI YANG MUDAH
Pilihlah sebuah warna, campurkan tiga warna, pilih preset animasi yang menyenangkan, atau buat sendiri animasimu. Pilihannya ada di tanganmu! Kamu bahkan bisa menyinkronkan mouse-mu dengan Logitech G LIGHTSYNC gear lainnya untuk kombinasi luar biasa.
AUDIO VISUALIZER
Mainkan musik, film, game, bahkan setiap audio, maka G203 akan menghadirkan warna sesuai irama musik.
SENSOR KELAS GAMING
Dapatkan penelusuran kursor yang akurat dan kinerja yang responsif berkat sensor kelas gaming. Dengan sensitivitas yang dapat disesuaikan antara 200-8.000 DPI, pilih level yang tepat sesuai dengan preferensi bermainmu. Gunakan software Logitech G HUB untuk memprogram hingga maksimal 5 preset.
PENGENCANGAN TOMBOL YANG DIOPTIMALKAN
Tombol kiri dan kanan utama memiliki sistem pengencangan tombol pegas logam eksklusif Logitech G yang menghadirkan aktuasi tombol yang akurat dan pengalaman yang konsisten—klik demi klik.
'@
$ws.Range("B3").Value = $desc1
$ws.Range("C3").Value = 269000
$ws.Range("D3").Value = 4.9
$ws.Range("E3").Value = 'https://www.tokopedia.com/logitech-g/logitech-g203-mouse-gaming-wired-rgb-lightsync-with-macro-fs-hitam?extParam=ivf%3Dtrue%26keyword%3Dmouse&src=topads&t_id=1747698251372&t_st=3&t_pp=search_result&t_efo=search_pure_goods_card&t_ef=goods_search&t_sm=&t_spt=search_result'
$ws.Rows.Item(3).AutoFit() | Out-Null

# Row 4: Aula F75 Mechanical Keyboard
$ws.Range("A4").Value = 'Aula F75 Mechanical Keyboard'
$desc2 = @'
F75 KEYBOARD AULA
Produk : Mechanical Keyboard
Brand : Aula
Model : F75
Bahan : Plastic
Plate : PC
Tipe : Keyboard Set
Lights : SMD LED RGB Backlit - 16.8 Million Color
LED Direction : South-facing(Side-printed) / North-facing(others)
Interface : USB Type-C + Wireless 2.4 Ghz + Bluetooth
Support : Windows + Mac + iOs + Android
Software : Aula Software
Berat : 1.1 kg
Size : 322.7 x 143.2 x 43.1 ± 1 mm
'@
$ws.Range("B4").Value = $desc2
$ws.Range("C4").Value = 739000
$ws.Range("D4").Value = 4.9
$ws.Range("E4").Value = 'https://tk.tokopedia.com/ZShpJYaYx/'
$ws.Rows.Item(4).AutoFit() | Out-Null

# Row 5: Terrel Sportswear Basic Tee
$ws.Range("A5").Value = 'Terrel Sportswear Basic Tee White Tshirt Baju Kaos Olah Raga Dry Fit Lari Running Gym Pria'
$desc3 = @'
Bahan: Polyester 
Baju Polyester menggunakan bahan berkualitas tinggi, baju ini memberikan kenyamanan yang memiliki sifat tahan lama, sehingga baju ini akan tetap terlihat bagus bahkan setelah banyak pemakaian.
Fitur :
1. Quick Dry: Baju ini dapat mengeringkan keringat dengan cepat. Anda akan tetap merasa segar dan nyaman bahkan saat beraktivitas fisik dengan intensitas tinggi.
2. Light Weight: Dengan bahan polyester yang ringan, baju ini memberikan kenyamanan saat digunakan. Bobotnya yang minimalis memungkinkan Anda untuk bergerak dengan bebas tanpa merasa terbebani oleh pakaian.
3. Anti Wrinkle: Sifat anti wrinkle yang membuatnya bebas dari kerutan. Anda dapat beraktivitas sepanjang hari tanpa khawatir harus menyeterika pakaian berulang kali.
Do & dont :
1. Jangan di cuci mesin.
2. Jangan di setrika
'@
$ws.Range("B5").Value = $desc3
$ws.Range("C5").Value = 72700
$ws.Range("D5").Value = 4.9
$ws.Range("E5").Value = 'https://shopee.co.id/Terrel-Sportswear-Basic-Tee-White-Tshirt-Baju-Kaos-Olah-Raga-Dry-Fit-Lari-Running-Gym-Pria-i.131221669.8935434906'
$ws.Rows.Item(5).AutoFit() | Out-Null

$ws.Columns.Item(1).ColumnWidth = 27.666666666666668

$ws.Range("D29").Select()

Write-Host "done"
